# ============================================================================
# Edit script: restructure PlayerPerformance_4224.xlsx
#   - Insert new "Player Info" sheet (first position)
#   - Keep "ODI Batting" / "ODI Bowling" (rename MATCH_CARD_LINK -> MATCH_CODE,
#     replace the full scorecard URL with just the numeric match code, and
#     drop now-empty INNING_NUMBER placeholder cells on "did not bat" rows)
#   - Append new "ODI Batting Extra" sheet (last position)
# ============================================================================

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

function Write-RowValues($ws, [int]$r, $values) {
    # $values is a flat, 1-D array of scalars for row $r starting at column 1.
    # $null   -> clear the cell (matches source: no data)
    # [int]   -> written as a genuine number
    # other   -> forced to literal text (matches source: every data cell in
    #            the original workbook is a plain inline string, never an
    #            auto-detected number/date)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($r, $col)
        $v = $values[$i]
        if ($null -eq $v) {
            $cell.ClearContents()
        } elseif ($v -is [int]) {
            $cell.NumberFormat = "General"
            $cell.Value = $v
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = [string]$v
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---- Restructure worksheets -----------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$wsPlayerInfo.Name = "Player Info"

$wsExtra = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsExtra.Name = "ODI Batting Extra"

$wsBatting = $wb.Worksheets("ODI Batting")
$wsBowling = $wb.Worksheets("ODI Bowling")

# ---- Player Info sheet --------------------------------------------------
Write-RowValues $wsPlayerInfo 1 @('ID', 'NAME', 'BATTING_HAND', 'BOWL_STYLE')
Set-HeaderStyle ($wsPlayerInfo.Range("A1:D1"))
Write-RowValues $wsPlayerInfo 2 @('4224', 'Ashton C Agar', 'Left Handed', 'Left Arm Orthodox')
$wsPlayerInfo.Range("A1").Select()

# ---- ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE -------------------
$wsBatting.Range("D1").NumberFormat = "@"
$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @('3832', '3834', '4069', '4071', '4166', '4167', '4168', '4169', '4170', '4398', '4399', '4400', '4419', '4437', '4486', '4594', '4645', '4646', '4660', '4663', '4732')
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $r = $i + 2
    $cell = $wsBatting.Cells.Item($r, 4)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$battingCodes[$i]
}

$battingRowsClearB = @(3, 11, 18, 20)
foreach ($rr in $battingRowsClearB) {
    $wsBatting.Cells.Item($rr, 2).ClearContents()
}

# ---- ODI Bowling sheet: MATCH_CARD_LINK -> MATCH_CODE -------------------
$wsBowling.Range("B1").NumberFormat = "@"
$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @('3832', '3834', '4069', '4071', '4166', '4167', '4168', '4169', '4170', '4398', '4399', '4400', '4419', '4437', '4486', '4594', '4645', '4646', '4660', '4663', '4732')
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $r = $i + 2
    $cell = $wsBowling.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$bowlingCodes[$i]
}

# ---- ODI Batting Extra sheet (new) ---------------------------------------
Write-RowValues $wsExtra 1 @('MATCH_CODE', 'BATTING_POSITION', 'NUM_4', 'NUM_6', 'PERCENT_RUNS_OF_TOTAL', 'MAN_OF_MATCH')
Set-HeaderStyle ($wsExtra.Range("A1:F1"))

Write-RowValues $wsExtra 2 @('3834', 8, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 3 @('4069', 8, '0', '0', $null, 'NO')
Write-RowValues $wsExtra 4 @('4071', 8, '1', '0', '3.07%', 'NO')
Write-RowValues $wsExtra 5 @('4166', 7, '4', '0', '18.69%', 'NO')
Write-RowValues $wsExtra 6 @('4167', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 7 @('4168', 8, '2', '1', '10.46%', 'NO')
Write-RowValues $wsExtra 8 @('4169', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 9 @('4170', 8, '0', '0', $null, 'NO')
Write-RowValues $wsExtra 10 @('4398', 7, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 11 @('4399', 7, '3', '0', '8.22%', 'NO')
Write-RowValues $wsExtra 12 @('4400', 8, '0', '0', '3.85%', 'NO')
Write-RowValues $wsExtra 13 @('4419', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 14 @('4437', 8, '2', '0', '9.69%', 'NO')
Write-RowValues $wsExtra 15 @('4486', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 16 @('4594', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 17 @('4645', 8, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 18 @('4646', 8, '0', '0', $null, 'NO')
Write-RowValues $wsExtra 19 @('4660', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 20 @('4663', $null, $null, $null, $null, 'NO')
Write-RowValues $wsExtra 21 @('4732', 9, '0', '1', '6.32%', 'NO')

$wsExtra.Range("A1").Select()

# ---- Final touches ----------------------------------------------------
$wb.Worksheets.Item(1).Select()
